$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — column F updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3120
$wsExhibit.Range("F5").Value = 6835
$wsExhibit.Range("F6").Value = 1895
$wsExhibit.Range("F7").Value = 5
$wsExhibit.Range("F8").Value = 61
$wsExhibit.Range("F14").Value = 161

# Sheet "全部类型" (sheet4) — column F updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3120
$wsAll.Range("F6").Value = 6835
$wsAll.Range("F7").Value = 1895
$wsAll.Range("F8").Value = 5
$wsAll.Range("F9").Value = 61
$wsAll.Range("F15").Value = 161
